# Apply scheduled-runner price/profit updates across all Cerberus sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 91.59999999999999
$ws.Range("I5").Value = 99.625
$ws.Range("K5").Value = 99.625
$ws.Range("M5").Value = 15.375

$ws.Range("H17").Value = 552.14703
$ws.Range("J17").Value = 552.14703
$ws.Range("L17").Value = 1656.44109
$ws.Range("N17").Value = -1992.44109

$ws.Range("H31").Value = 4118.1113
$ws.Range("I31").Value = 4008.625
$ws.Range("K31").Value = 12025.875
$ws.Range("M31").Value = -11795.875

$ws.Range("H33").Value = 339.57144
$ws.Range("I33").Value = 350.33334
$ws.Range("J33").Value = 331.5
$ws.Range("K33").Value = 350.33334
$ws.Range("L33").Value = 331.5
$ws.Range("M33").Value = -121.33334
$ws.Range("N33").Value = -789.5

$ws.Range("H76").Value = 500
$ws.Range("I76").Value = 500
$ws.Range("K76").Value = 500
$ws.Range("M76").Value = -185

$ws.Range("H79").Value = 500
$ws.Range("I79").Value = 500
$ws.Range("K79").Value = 500
$ws.Range("M79").Value = 592

$ws.Range("H100").Value = 2826.5
$ws.Range("J100").Value = 2499.5
$ws.Range("L100").Value = 2499.5
$ws.Range("N100").Value = -3581.5

$ws.Range("H103").Value = 413.25
$ws.Range("I103").Value = 351.75
$ws.Range("K103").Value = 1055.25
$ws.Range("M103").Value = -469.25


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 54996
$ws.Range("J37").Value = 54992
$ws.Range("L37").Value = 54992
$ws.Range("N37").Value = -55538

$ws.Range("H61").Value = 5233.1304
$ws.Range("I61").Value = 4038.111
$ws.Range("J61").Value = 9535.200000000001
$ws.Range("K61").Value = 4038.111
$ws.Range("L61").Value = 9535.200000000001
$ws.Range("M61").Value = -3826.111
$ws.Range("N61").Value = -9959.200000000001

$ws.Range("H63").Value = 3135.5
$ws.Range("I63").Value = 3085
$ws.Range("K63").Value = 3085
$ws.Range("M63").Value = -2399

$ws.Range("H66").Value = 3135.5
$ws.Range("I66").Value = 3085
$ws.Range("K66").Value = 15425
$ws.Range("M66").Value = -11993

$ws.Range("H88").Value = 12532.728
$ws.Range("I88").Value = 1898.25
$ws.Range("J88").Value = 18609.572
$ws.Range("K88").Value = 1898.25
$ws.Range("L88").Value = 18609.572
$ws.Range("M88").Value = -1492.25
$ws.Range("N88").Value = -19421.572

$ws.Range("H91").Value = 12532.728
$ws.Range("I91").Value = 1898.25
$ws.Range("J91").Value = 18609.572
$ws.Range("K91").Value = 1898.25
$ws.Range("L91").Value = 18609.572
$ws.Range("M91").Value = -494.25
$ws.Range("N91").Value = -21417.572

$ws.Range("H136").Value = 5233.1304
$ws.Range("I136").Value = 4038.111
$ws.Range("J136").Value = 9535.200000000001
$ws.Range("K136").Value = 12114.333
$ws.Range("L136").Value = 28605.6
$ws.Range("M136").Value = -9564.332999999999
$ws.Range("N136").Value = -33705.60000000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 237.38461
$ws.Range("I22").Value = 232.6
$ws.Range("K22").Value = 232.6
$ws.Range("M22").Value = -59.59999999999999

$ws.Range("H134").Value = 8721.932000000001
$ws.Range("I134").Value = 7928.0967
$ws.Range("K134").Value = 23784.2901
$ws.Range("M134").Value = -21249.2901


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 63.57143
$ws.Range("J7").Value = 83
$ws.Range("L7").Value = 83
$ws.Range("N7").Value = -309

$ws.Range("H16").Value = 2509.0588
$ws.Range("J16").Value = 4142.625
$ws.Range("L16").Value = 4142.625
$ws.Range("N16").Value = -4716.625

$ws.Range("H22").Value = 331.14285
$ws.Range("I22").Value = 285.42856
$ws.Range("K22").Value = 285.42856
$ws.Range("M22").Value = 64.57144

$ws.Range("H94").Value = 1281.2667
$ws.Range("I94").Value = 833
$ws.Range("J94").Value = 1673.5
$ws.Range("K94").Value = 833
$ws.Range("L94").Value = 1673.5
$ws.Range("M94").Value = -382
$ws.Range("N94").Value = -2575.5

$ws.Range("H113").Value = 2509.0588
$ws.Range("J113").Value = 4142.625
$ws.Range("L113").Value = 4142.625
$ws.Range("N113").Value = -8482.625


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1999.6666
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H128").Value = 125000
$ws.Range("I128").Value = 125000
$ws.Range("K128").Value = 375000
$ws.Range("M128").Value = -370020

$ws.Range("H129").Value = 3328.3333
$ws.Range("J129").Value = 6099.75
$ws.Range("L129").Value = 18299.25
$ws.Range("N129").Value = -28299.25


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3117.55
$ws.Range("I122").Value = 2353.3076
$ws.Range("J122").Value = 4536.857
$ws.Range("K122").Value = 7059.9228
$ws.Range("L122").Value = 13610.571
$ws.Range("M122").Value = -4609.9228
$ws.Range("N122").Value = -18510.571

$ws.Range("H134").Value = 50829.4
$ws.Range("J134").Value = 50829.4
$ws.Range("L134").Value = 152488.2
$ws.Range("N134").Value = -157558.2


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1917.7407
$ws.Range("I61").Value = 1917.7407
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1917.7407
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1715.7407
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 1917.7407
$ws.Range("I113").Value = 1917.7407
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1917.7407
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 252.2592999999999
$ws.Range("N113").ClearContents()


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1971.7646
$ws.Range("I100").Value = 1097.2727
$ws.Range("K100").Value = 2194.5454
$ws.Range("M100").Value = -1653.5454

$ws.Range("H107").Value = 1738.6
$ws.Range("I107").Value = 812.1429000000001
$ws.Range("J107").Value = 3900.3333
$ws.Range("K107").Value = 2436.4287
$ws.Range("L107").Value = 11700.9999
$ws.Range("M107").Value = -516.4287000000004
$ws.Range("N107").Value = -15540.9999

$ws.Range("H122").Value = 2587.5789
$ws.Range("I122").Value = 2559.1667
$ws.Range("J122").Value = 3099
$ws.Range("K122").Value = 7677.500100000001
$ws.Range("L122").Value = 9297
$ws.Range("M122").Value = -5227.500100000001
$ws.Range("N122").Value = -14197

$ws.Range("H132").Value = 20003424
$ws.Range("I132").Value = 22225204
$ws.Range("K132").Value = 66675612
$ws.Range("M132").Value = -66673082

$ws.Range("H136").Value = 5477.643
$ws.Range("I136").Value = 4931.4634
$ws.Range("J136").Value = 6970.533
$ws.Range("K136").Value = 14794.3902
$ws.Range("L136").Value = 20911.599
$ws.Range("M136").Value = -12244.3902
$ws.Range("N136").Value = -26011.599

